$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of link data (row 5 is left blank, matching the source sheet's gap before row 6)
$ws.Range("A6").Value = "Combinational Circuits"
$ws.Range("B6").Value = "NOR (Half Add/Sub)"
$ws.Range("C6").Value = "https://www.youtube.com/watch?v=P_UW41wMvpM"

# Give the new link cell the same visual "Hyperlink" look used by the other
# link cells, without actually inserting a clickable hyperlink.
$ws.Range("C6").Style = "Hyperlink"

# Columns got a bit wider to fit the new, longer text.
$ws.Columns.Item(1).ColumnWidth = 18.6
$ws.Columns.Item(2).ColumnWidth = 16.1

# Restore the selection to where the author last left the cursor.
[void]$ws.Range("C13").Select()
